# Weekly data refresh: a new weekly price observation for Brócoli at the
# "Feria Lagunitas de Puerto Montt" market is inserted as a new record
# (row 371), pushing the previously existing records (old rows 371-458)
# down by one row (to 372-459).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 371; this shifts rows 371:458 down to 372:459
# and grows the sheet from A1:R458 to A1:R459.
$ws.Rows(371).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(371, 1).Value2 = 4
$ws.Cells.Item(371, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(371, 3).Value2 = "Los Lagos"
$ws.Cells.Item(371, 4).Value2 = 44932
$ws.Cells.Item(371, 5).Value2 = 10
$ws.Cells.Item(371, 6).Value2 = 100112023
$ws.Cells.Item(371, 7).Value2 = "Brócoli"
$ws.Cells.Item(371, 8).Value2 = "Sin especificar"
$ws.Cells.Item(371, 9).Value2 = "Primera"
$ws.Cells.Item(371, 10).Value2 = 1400
$ws.Cells.Item(371, 11).Value2 = 1500
$ws.Cells.Item(371, 12).Value2 = 1500
$ws.Cells.Item(371, 13).Value2 = 1500
$ws.Cells.Item(371, 14).Value2 = "`$/unidad"
$ws.Cells.Item(371, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(371, 16).Value2 = 1500
$ws.Cells.Item(371, 17).Value2 = 1
$ws.Cells.Item(371, 18).Value2 = "Hortaliza"
